$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Columns D (Price) and E (Volume(1h)) hold text-like values (e.g. "1.009", "20.519.66",
# "  +1.63%  ") that must remain literal text rather than being auto-converted by Excel
# into numbers/scientific notation, so we force the cell NumberFormat to Text ("@") first.

# Row 2
$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "20.519.66"
$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = "  +1.63%  "

# Row 3
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "1.473.66"
$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = "  +2.37%  "

# Row 4
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = "1.009"
$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = "  +0.09%  "

# Row 5
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "0.9538"
$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = "  +2.97%  "

# Row 6
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "277.44"
$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = "  +0.31%  "

# Row 7
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.3609"
$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = "  -1.38%  "

# Row 8
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.3055"
$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = "  -2.25%  "

# Row 9
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "39.39"
$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = "  +0.20%  "

# Row 10
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "1.056"
$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = "  +3.72%  "

# Row 11
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.06636"
$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = "  +1.91%  "

# Row 12
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "1.003"
$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = "  +0.06%  "

# Row 13
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "5.506"
$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = "  +1.93%  "

# Row 14
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "18.08"
$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = "  +3.08%  "

# Row 15
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "6.175"
$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = "  +1.59%  "

# Row 16
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "0.9565"
$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = "  +2.10%  "

# Row 17
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "0.00001027"
$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = "  +1.27%  "

# Row 18
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "1.475.54"
$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = "  +1.90%  "

# Row 19
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "0.05944"
$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = "  +6.19%  "

# Row 20
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "69.05"
$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = "  +2.74%  "

# Row 21
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "5.491"
$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = "  +1.35%  "

# Row 22
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "14.47"
$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = "  +0.40%  "

# Row 23
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "11.14"
$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = "  +2.64%  "

# Row 24
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "2.251"
$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = "  -1.02%  "

# Row 25
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "20.564.64"
$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = "  +1.37%  "

# Row 26
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "143.13"
$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = "  +5.56%  "

# Row 27
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "2.129"
$r = $ws.Range("E27")
$r.NumberFormat = "@"
$r.Value = "  -2.79%  "

# Row 28
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "17.14"
$r = $ws.Range("E28")
$r.NumberFormat = "@"
$r.Value = "  +0.97%  "

# Row 29
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "1.636.79"
$r = $ws.Range("E29")
$r.NumberFormat = "@"
$r.Value = "  +2.24%  "

# Row 30
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "113.54"
$r = $ws.Range("E30")
$r.NumberFormat = "@"
$r.Value = "  +2.83%  "

# Row 31
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "3.940"

# Row 32
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "4.997"
$r = $ws.Range("E32")
$r.NumberFormat = "@"
$r.Value = "  +2.75%  "

# Row 33
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "0.07983"
$r = $ws.Range("E33")
$r.NumberFormat = "@"
$r.Value = "  +4.55%  "

# Row 34
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "0.8063"
$r = $ws.Range("E34")
$r.NumberFormat = "@"
$r.Value = "  -0.84%  "

# Row 35
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "1.514"
$r = $ws.Range("E35")
$r.NumberFormat = "@"
$r.Value = "  +0.19%  "

# Row 36
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "1.215"
$r = $ws.Range("E36")
$r.NumberFormat = "@"
$r.Value = "  +7.00%  "

# Row 37
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "0.05835"
$r = $ws.Range("E37")
$r.NumberFormat = "@"
$r.Value = "  -2.16%  "

# Row 38
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "4.708"
$r = $ws.Range("E38")
$r.NumberFormat = "@"
$r.Value = "  +0.25%  "

# Row 39
$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = "  +3.10%  "

# Row 40
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.9567"
$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = "  +2.27%  "

# Row 41
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "10.33"
$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = "  +0.71%  "

# Row 42
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "0.1875"
$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = "  +2.62%  "

# Row 43
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "7.430"
$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = "  +4.71%  "

# Row 44
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.5280"
$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = "  +0.83%  "

# Row 45
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "3.520"
$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = "  +0.11%  "

# Row 46
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "12.18"
$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = "  +0.84%  "

# Row 47
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "118.17"
$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = "  -1.51%  "

# Row 48
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "0.5186"
$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = "  +0.59%  "

# Row 49
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "1.810"
$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = "  +2.45%  "

# Row 50
$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = "  +2.00%  "

# Row 51
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "0.9786"
$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = "  -1.10%  "
